# Add a new "Serviced by " column (O) to the Card17 sheet and fix up the
# "Correction" header / the now-populated N column "nan" placeholders,
# mirroring the layout already used on the sibling Card* sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card17")

# Clone the header style (bold, bordered, centered - style index 1 on the
# other sheets) from N1 onto the brand-new O1 cell, then fix up the text of
# both cells.
$ws.Range("N1").Copy($ws.Range("O1")) | Out-Null

$ws.Range("N1").Value = "Correction"
$ws.Range("O1").Value = "Serviced by "

# Rows 2-12: N was an empty placeholder cell before; it now holds the same
# "nan" placeholder text used throughout the rest of the row. O is the new
# column and stays blank, same as on the other already-migrated sheets.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 14).Value = "nan"
}
